$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.511549711227417
$ws.Range("B1").Value = 2.018351793289185
$ws.Range("C1").Value = 3.41135835647583
$ws.Range("D1").Value = 1.366315603256226
$ws.Range("E1").Value = 0.8139936327934265
